$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row 12 before the existing "row 15" content, pushing nothing
# else down (rows 12-14 were already blank in the sheet, so we just author
# the row 12 cells directly).
$ws.Cells.Item(12, 2).Value = "https://hildemuz.itch.io/banana-man"
$ws.Cells.Item(12, 1).Value = "Banana Man"
$ws.Cells.Item(12, 3).Value = "You can use it in your games freely if you download from here and provide information about me in your game."

# Row height for the new row
$ws.Rows.Item(12).RowHeight = 24

# New font/style applied to C12: Arial 19pt, color FF505050
$c12 = $ws.Cells.Item(12, 3)
$c12.Font.Name = "Arial"
$c12.Font.Size = 19
$c12.Font.Color = 5263440

# Update the view: top-left cell and selection to match the new content
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C12").Select()
